$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.663.37"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "3.329.76"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.84"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.37"
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("D9").Value = "3.325.56"
$ws.Range("E9").Value = "  +2.29%  "
$ws.Range("E10").Value = "  +6.66%  "
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.73"
$ws.Range("E12").Value = "  +4.87%  "
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "690.96"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "3.876.03"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.41"
$ws.Range("E16").Value = "  +2.70%  "
$ws.Range("D17").Value = "67.665.17"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "3.334.37"
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("E20").Value = "  +2.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.02"
$ws.Range("E21").Value = "  +4.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.893"
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("E23").Value = "  +5.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.91"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.78"
$ws.Range("E25").Value = "  +4.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.90"
$ws.Range("E27").Value = "  +2.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.45"
$ws.Range("E28").Value = "  +6.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.90"
$ws.Range("E29").Value = "  +2.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.53"
$ws.Range("E30").Value = "  +3.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.03"
$ws.Range("E31").Value = "  +6.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "569.22"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.99"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("E34").Value = "  +3.36%  "
$ws.Range("D35").Value = "3.714.07"
$ws.Range("E35").Value = "  -1.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "57.16"
$ws.Range("E37").Value = "  +3.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.27"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.10"
$ws.Range("E39").Value = "  +12.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.133"
$ws.Range("E40").Value = "  +4.77%  "
$ws.Range("E41").Value = "  +7.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.61"
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.34"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("D44").Value = "0.0₃0671"
$ws.Range("E44").Value = "  +3.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.334"
$ws.Range("E45").Value = "  +3.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0407"
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("E47").Value = "  +6.28%  "
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  +1.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.09"
$ws.Range("E51").Value = "  +4.38%  "
